# Updated cryptos list on Fri Aug  4 04:59:19 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (col D) / "Volume(1h)" (col E) figures scraped from
# coinranking.com, and corrects the ranking order of FraxShare and
# XinFinNetwork (rows 42 & 43 swapped places).
#
# NOTE: column D prices are stored as plain text (coinranking formats them
# with "." as a thousands separator AND a decimal point, e.g. "29.159.95",
# which is not a legal number), so a leading "'" is used for values that
# would otherwise look like an ordinary decimal number (e.g. "241.63") to
# stop Excel's COM layer from silently converting them to a numeric type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 42/43: FraxShare and XinFinNetwork swapped positions -------------
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'5.902"
$ws.Range("E42").Value = "  -2.65%  "

$ws.Range("B43").Value = "XinFinNetwork"
$ws.Range("C43").Value = "https://coinranking.com/coin/77jGXSqWJ1ofG+xinfinnetwork-xdc"
$ws.Range("D43").Value = "'0.08486"
$ws.Range("E43").Value = "  +0.45%  "

# --- Price / Volume(1h) refresh for every other row ------------------------
$ws.Range("D2").Value = "29.159.95"
$ws.Range("E2").Value = "  +0.15%  "

$ws.Range("D3").Value = "1.833.03"
$ws.Range("E3").Value = "  -0.16%  "

$ws.Range("D5").Value = "'241.63"
$ws.Range("E5").Value = "  +0.79%  "

$ws.Range("E6").Value = "  -2.32%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.31%  "

$ws.Range("D9").Value = "'0.2940"
$ws.Range("E9").Value = "  -1.57%  "

$ws.Range("D10").Value = "'22.81"
$ws.Range("E10").Value = "  -1.42%  "

$ws.Range("D11").Value = "'0.07752"
$ws.Range("E11").Value = "  +1.36%  "

$ws.Range("D12").Value = "1.835.72"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").Value = "'4.985"
$ws.Range("E13").Value = "  -0.75%  "

$ws.Range("D14").Value = "'0.6687"
$ws.Range("E14").Value = "  -1.41%  "

$ws.Range("D15").Value = "'82.89"
$ws.Range("E15").Value = "  -4.57%  "

$ws.Range("D16").Value = "'6.097"
$ws.Range("E16").Value = "  -0.95%  "

$ws.Range("D17").Value = "'0.000008355"
$ws.Range("E17").Value = "  +1.51%  "

$ws.Range("D18").Value = "29.173.71"
$ws.Range("E18").Value = "  +0.30%  "

$ws.Range("D19").Value = "2.086.21"
$ws.Range("E19").Value = "  +0.89%  "

$ws.Range("D20").Value = "'228.71"
$ws.Range("E20").Value = "  -0.15%  "

$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +0.12%  "

$ws.Range("D23").Value = "'7.164"
$ws.Range("E23").Value = "  -2.25%  "

$ws.Range("D24").Value = "'0.9997"
$ws.Range("E24").Value = "  -0.04%  "

$ws.Range("D25").Value = "'159.45"
$ws.Range("E25").Value = "  -1.15%  "

$ws.Range("D26").Value = "'0.1407"
$ws.Range("E26").Value = "  -1.72%  "

$ws.Range("D28").Value = "'18.00"
$ws.Range("E28").Value = "  -0.16%  "

$ws.Range("D29").Value = "'1.511"
$ws.Range("E29").Value = "  +0.88%  "

$ws.Range("E30").Value = "  -3.17%  "

$ws.Range("D31").Value = "'4.039"
$ws.Range("E31").Value = "  -2.42%  "

$ws.Range("D32").Value = "'1.190"
$ws.Range("E32").Value = "  +0.12%  "

$ws.Range("E33").Value = "  -0.41%  "

$ws.Range("D34").Value = "'1.869"
$ws.Range("E34").Value = "  +1.26%  "

$ws.Range("D35").Value = "'0.7485"
$ws.Range("E35").Value = "  -0.72%  "

$ws.Range("D36").Value = "'1.138"
$ws.Range("E36").Value = "  +0.62%  "

$ws.Range("D37").Value = "'2.640"
$ws.Range("E37").Value = "  -1.61%  "

$ws.Range("D38").Value = "1.274.05"
$ws.Range("E38").Value = "  -2.98%  "

$ws.Range("D39").Value = "'0.01796"
$ws.Range("E39").Value = "  -1.37%  "

$ws.Range("D40").Value = "'2.736"
$ws.Range("E40").Value = "  +0.86%  "

$ws.Range("D41").Value = "'0.9288"
$ws.Range("E41").Value = "  -0.92%  "

$ws.Range("D45").Value = "'101.99"
$ws.Range("E45").Value = "  -2.78%  "

$ws.Range("D46").Value = "1.981.82"
$ws.Range("E46").Value = "  +0.27%  "

$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  -1.59%  "

$ws.Range("D50").Value = "'63.03"
$ws.Range("E50").Value = "  -1.35%  "

$ws.Range("E51").Value = "  -0.84%  "
